{"js": "// Replace the date line and every \"A\u00d7B=C\" answer cell with its updated value.\n// Every search string below is unique within the document, so a plain\n// body.search + insertText(replace) round-trip is unambiguous.\nconst replacements = [\n  [\"2024-10-17 Thursday\", \"2024-10-18 Friday\"],\n  [\"571\u00d79=5139\", \"778\u00d73=2334\"],\n  [\"355\u00d77=2485\", \"406\u00d79=3654\"],\n  [\"693\u00d77=4851\", \"518\u00d79=4662\"],\n  [\"829\u00d77=5803\", \"567\u00d77=3969\"],\n  [\"726\u00d77=5082\", \"166\u00d79=1494\"],\n  [\"551\u00d74=2204\", \"859\u00d76=5154\"],\n  [\"394\u00d76=2364\", \"283\u00d73=849\"],\n  [\"367\u00d75=1835\", \"338\u00d77=2366\"],\n  [\"562\u00d73=1686\", \"567\u00d73=1701\"],\n  [\"899\u00d72=1798\", \"556\u00d73=1668\"],\n  [\"241\u00d77=1687\", \"965\u00d78=7720\"],\n  [\"914\u00d72=1828\", \"247\u00d73=741\"],\n  [\"161\u00d76=966\", \"947\u00d74=3788\"],\n  [\"642\u00d73=1926\", \"123\u00d78=984\"],\n  [\"714\u00d77=4998\", \"410\u00d73=1230\"],\n  [\"560\u00d77=3920\", \"643\u00d76=3858\"],\n  [\"231\u00d73=693\", \"460\u00d76=2760\"],\n  [\"898\u00d79=8082\", \"101\u00d77=707\"],\n  [\"342\u00d75=1710\", \"656\u00d77=4592\"],\n  [\"574\u00d76=3444\", \"713\u00d79=6417\"],\n  [\"991\u00d73=2973\", \"481\u00d73=1443\"],\n  [\"969\u00d78=7752\", \"649\u00d74=2596\"],\n  [\"980\u00d79=8820\", \"831\u00d78=6648\"],\n  [\"188\u00d73=564\", \"247\u00d74=988\"],\n  [\"602\u00d76=3612\", \"850\u00d76=5100\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" answer cell with its updated value.\n# Each search string is unique within the document, so Find/Replace with\n# wdReplaceAll (2) on each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-10-17 Thursday\", \"2024-10-18 Friday\"),\n  @(\"571\u00d79=5139\", \"778\u00d73=2334\"),\n  @(\"355\u00d77=2485\", \"406\u00d79=3654\"),\n  @(\"693\u00d77=4851\", \"518\u00d79=4662\"),\n  @(\"829\u00d77=5803\", \"567\u00d77=3969\"),\n  @(\"726\u00d77=5082\", \"166\u00d79=1494\"),\n  @(\"551\u00d74=2204\", \"859\u00d76=5154\"),\n  @(\"394\u00d76=2364\", \"283\u00d73=849\"),\n  @(\"367\u00d75=1835\", \"338\u00d77=2366\"),\n  @(\"562\u00d73=1686\", \"567\u00d73=1701\"),\n  @(\"899\u00d72=1798\", \"556\u00d73=1668\"),\n  @(\"241\u00d77=1687\", \"965\u00d78=7720\"),\n  @(\"914\u00d72=1828\", \"247\u00d73=741\"),\n  @(\"161\u00d76=966\", \"947\u00d74=3788\"),\n  @(\"642\u00d73=1926\", \"123\u00d78=984\"),\n  @(\"714\u00d77=4998\", \"410\u00d73=1230\"),\n  @(\"560\u00d77=3920\", \"643\u00d76=3858\"),\n  @(\"231\u00d73=693\", \"460\u00d76=2760\"),\n  @(\"898\u00d79=8082\", \"101\u00d77=707\"),\n  @(\"342\u00d75=1710\", \"656\u00d77=4592\"),\n  @(\"574\u00d76=3444\", \"713\u00d79=6417\"),\n  @(\"991\u00d73=2973\", \"481\u00d73=1443\"),\n  @(\"969\u00d78=7752\", \"649\u00d74=2596\"),\n  @(\"980\u00d79=8820\", \"831\u00d78=6648\"),\n  @(\"188\u00d73=564\", \"247\u00d74=988\"),\n  @(\"602\u00d76=3612\", \"850\u00d76=5100\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  # wdFindContinue = 1, wdReplaceAll = 2\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
